$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 97.65533411488863
$ws.Range("D3").Value = 97.57673667205169
$ws.Range("D4").Value = 96.44153957879448
$ws.Range("D5").Value = 95.99198396793587
$ws.Range("D6").Value = 95.37815126050421
$ws.Range("D7").Value = 94.67680608365019
$ws.Range("D8").Value = 90.20332717190388
$ws.Range("D9").Value = 89.1498039836851
$ws.Range("D10").Value = 92.31096399578379
